$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2749.5
$ws.Range("I18").Value = 2749.5
$ws.Range("K18").Value = 2749.5
$ws.Range("M18").Value = -2465.5
$ws.Range("H32").Value = 2527.3333
$ws.Range("J32").Value = 2678
$ws.Range("L32").Value = 2678
$ws.Range("N32").Value = -3330
$ws.Range("H53").Value = 328.1111
$ws.Range("I53").Value = 196.33333
$ws.Range("J53").Value = 394
$ws.Range("K53").Value = 196.33333
$ws.Range("L53").Value = 394
$ws.Range("M53").Value = 440.66667
$ws.Range("N53").Value = -1668
$ws.Range("H64").Value = 5463.6665
$ws.Range("I64").Value = 5399
$ws.Range("J64").Value = 5496
$ws.Range("K64").Value = 5399
$ws.Range("L64").Value = 5496
$ws.Range("M64").Value = -5151
$ws.Range("N64").Value = -5992
$ws.Range("H67").Value = 5463.6665
$ws.Range("I67").Value = 5399
$ws.Range("J67").Value = 5496
$ws.Range("K67").Value = 5399
$ws.Range("L67").Value = 5496
$ws.Range("M67").Value = -4541
$ws.Range("N67").Value = -7212
$ws.Range("H86").Value = 24685.572
$ws.Range("I86").Value = 24685.572
$ws.Range("K86").Value = 24685.572
$ws.Range("M86").Value = -23562.572
$ws.Range("H89").Value = 24685.572
$ws.Range("I89").Value = 24685.572
$ws.Range("K89").Value = 123427.86
$ws.Range("M89").Value = -117811.86
$ws.Range("H103").Value = 925.8
$ws.Range("I103").Value = 666
$ws.Range("J103").Value = 990.75
$ws.Range("K103").Value = 1998
$ws.Range("L103").Value = 2972.25
$ws.Range("M103").Value = -1412
$ws.Range("N103").Value = -4144.25
$ws.Range("H132").Value = 2222
$ws.Range("I132").Value = 2342.5
$ws.Range("K132").Value = 7027.5
$ws.Range("M132").Value = -4497.5
$ws.Range("H137").Value = 1172.5
$ws.Range("I137").Value = 1007
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 3021
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -471
$ws.Range("N137").Value = -11100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 99750.25
$ws.Range("I38").Value = 131667
$ws.Range("J38").Value = 4000
$ws.Range("K38").Value = 131667
$ws.Range("L38").Value = 4000
$ws.Range("M38").Value = -131200
$ws.Range("N38").Value = -4934
$ws.Range("H61").Value = 2407.1428
$ws.Range("I61").Value = 2308.3333
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2308.3333
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2096.3333
$ws.Range("N61").Value = -3424
$ws.Range("H88").Value = 2134.0667
$ws.Range("I88").Value = 1933.5
$ws.Range("J88").Value = 2535.2
$ws.Range("K88").Value = 1933.5
$ws.Range("L88").Value = 2535.2
$ws.Range("M88").Value = -1527.5
$ws.Range("N88").Value = -3347.2
$ws.Range("H91").Value = 2134.0667
$ws.Range("I91").Value = 1933.5
$ws.Range("J91").Value = 2535.2
$ws.Range("K91").Value = 1933.5
$ws.Range("L91").Value = 2535.2
$ws.Range("M91").Value = -529.5
$ws.Range("N91").Value = -5343.2
$ws.Range("H97").Value = 835.6667
$ws.Range("I97").Value = 848.2308
$ws.Range("K97").Value = 848.2308
$ws.Range("M97").Value = -352.2308
$ws.Range("H102").Value = 4000
$ws.Range("I102").Value = 4000
$ws.Range("K102").Value = 4000
$ws.Range("M102").Value = -2378
$ws.Range("H132").Value = 2597.8462
$ws.Range("I132").Value = 2697.5454
$ws.Range("J132").Value = 2049.5
$ws.Range("K132").Value = 8092.6362
$ws.Range("L132").Value = 6148.5
$ws.Range("M132").Value = -5562.6362
$ws.Range("N132").Value = -11208.5
$ws.Range("H136").Value = 2407.1428
$ws.Range("I136").Value = 2308.3333
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 6924.999899999999
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -4374.999899999999
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 139542
$ws.Range("J76").Value = 139542
$ws.Range("L76").Value = 139542
$ws.Range("N76").Value = -140172
$ws.Range("H79").Value = 139542
$ws.Range("J79").Value = 139542
$ws.Range("L79").Value = 139542
$ws.Range("N79").Value = -141726
$ws.Range("H99").Value = 1999.909
$ws.Range("I99").Value = 1999.9
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1999.9
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -501.9000000000001
$ws.Range("N99").Value = -4996
$ws.Range("H105").Value = 2461.9167
$ws.Range("I105").Value = 2482.6667
$ws.Range("J105").Value = 2399.6667
$ws.Range("K105").Value = 2482.6667
$ws.Range("L105").Value = 2399.6667
$ws.Range("M105").Value = -735.6667000000002
$ws.Range("N105").Value = -5893.6667
$ws.Range("H132").Value = 57889.668
$ws.Range("J132").Value = 57889.668
$ws.Range("L132").Value = 57889.668
$ws.Range("N132").Value = -68009.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = ""
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = ""
$ws.Range("H94").Value = 5283.6665
$ws.Range("I94").Value = 6061.75
$ws.Range("K94").Value = 6061.75
$ws.Range("M94").Value = -5610.75
$ws.Range("H99").Value = 3665.75
$ws.Range("J99").Value = 5331
$ws.Range("L99").Value = 5331
$ws.Range("N99").Value = -8327
$ws.Range("H126").Value = 3665.75
$ws.Range("J126").Value = 5331
$ws.Range("L126").Value = 15993
$ws.Range("N126").Value = -20933
$ws.Range("H141").Value = 184425.14
$ws.Range("J141").Value = 184425.14
$ws.Range("L141").Value = 184425.14
$ws.Range("N141").Value = -194785.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 224.5
$ws.Range("I2").Value = 69.666664
$ws.Range("J2").Value = 379.33334
$ws.Range("K2").Value = 417.999984
$ws.Range("L2").Value = 2276.00004
$ws.Range("M2").Value = -304.999984
$ws.Range("N2").Value = -2502.00004
$ws.Range("H32").Value = 3250
$ws.Range("J32").Value = 3250
$ws.Range("L32").Value = 9750
$ws.Range("N32").Value = -10316
$ws.Range("H38").Value = 2784.5
$ws.Range("I38").Value = 2784.5
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 8353.5
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -8006.5
$ws.Range("N38").Value = ""
$ws.Range("H92").Value = 478.83334
$ws.Range("I92").Value = 250.5
$ws.Range("K92").Value = 751.5
$ws.Range("M92").Value = 496.5
$ws.Range("H98").Value = 1458.3334
$ws.Range("I98").Value = 875
$ws.Range("J98").Value = 1750
$ws.Range("K98").Value = 2625
$ws.Range("L98").Value = 5250
$ws.Range("M98").Value = -1127
$ws.Range("N98").Value = -8246
$ws.Range("H109").Value = 1543.375
$ws.Range("I109").Value = 507.83334
$ws.Range("J109").Value = 4650
$ws.Range("K109").Value = 1523.50002
$ws.Range("L109").Value = 13950
$ws.Range("M109").Value = -483.5000199999999
$ws.Range("N109").Value = -16030
$ws.Range("H121").Value = 2926.6843
$ws.Range("I121").Value = 1809.6666
$ws.Range("J121").Value = 3136.125
$ws.Range("K121").Value = 5428.9998
$ws.Range("L121").Value = 9408.375
$ws.Range("M121").Value = -4118.9998
$ws.Range("N121").Value = -12028.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 14433
$ws.Range("I10").Value = 12000
$ws.Range("J10").Value = 15649.5
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 15649.5
$ws.Range("M10").Value = -11831
$ws.Range("N10").Value = -15987.5
$ws.Range("H62").Value = 90077
$ws.Range("I62").Value = 90077
$ws.Range("K62").Value = 90077
$ws.Range("M62").Value = -89391
$ws.Range("H65").Value = 90077
$ws.Range("I65").Value = 90077
$ws.Range("K65").Value = 270231
$ws.Range("M65").Value = -266799

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7844.8184
$ws.Range("I132").Value = 8129.2
$ws.Range("K132").Value = 24387.6
$ws.Range("M132").Value = -21857.6
$ws.Range("H136").Value = 3998
$ws.Range("J136").Value = 3998
$ws.Range("L136").Value = 11994
$ws.Range("N136").Value = -17094

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 85000
$ws.Range("J70").Value = 85000
$ws.Range("L70").Value = 85000
$ws.Range("N70").Value = -85630
$ws.Range("H73").Value = 85000
$ws.Range("J73").Value = 85000
$ws.Range("L73").Value = 85000
$ws.Range("N73").Value = -87184
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 75311.75
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 75311.75
$ws.Range("M75").Value = ""
$ws.Range("N75").Value = -77183.75
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 75311.75
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 225935.25
$ws.Range("M78").Value = ""
$ws.Range("N78").Value = -235295.25
$ws.Range("H100").Value = 25000724
$ws.Range("I100").Value = 33333966
$ws.Range("J100").Value = 993
$ws.Range("K100").Value = 66667932
$ws.Range("L100").Value = 1986
$ws.Range("M100").Value = -66667391
$ws.Range("N100").Value = -3068
$ws.Range("H121").Value = 55555
$ws.Range("J121").Value = 55555
$ws.Range("L121").Value = 55555
$ws.Range("N121").Value = -59049
$ws.Range("H122").Value = 2002.9333
$ws.Range("I122").Value = 1962
$ws.Range("J122").Value = 2166.6667
$ws.Range("K122").Value = 5886
$ws.Range("L122").Value = 6500.000100000001
$ws.Range("M122").Value = -3436
$ws.Range("N122").Value = -11400.0001
$ws.Range("H132").Value = 3064.8333
$ws.Range("I132").Value = 3064.8333
$ws.Range("K132").Value = 9194.499899999999
$ws.Range("M132").Value = -6664.499899999999
$ws.Range("H136").Value = 1947.5
$ws.Range("I136").Value = 1770.2667
$ws.Range("K136").Value = 5310.800099999999
$ws.Range("M136").Value = -2760.800099999999


Write-Output "Applied market price updates across all leve sheets."
